# ============================================================================
# PlayerPerformance_3903.xlsx edit script
#
# Summary of the change (per commit message / diff):
#  - Add a new "Player Info" sheet (first tab) with ID/NAME/BATTING_HAND/
#    BOWL_STYLE for player 3903 (Johnson Charles).
#  - Rename column D ("MATCH_CARD_LINK" -> "MATCH_CODE") on both the
#    "ODI Batting" and "ODI Bowling"/"ODI Bowling" sheets, and replace the
#    full howstat URL values with just the trailing MatchCode number
#    (kept as text).
#  - Add a new "ODI Batting Extra" sheet (last tab) with additional
#    per-match batting detail (batting position, 4s, 6s, % of team total,
#    man of the match).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Helper: turn a numeric-looking string into a forced-text literal so the
# COM layer stores it as a string cell (matching the source data, which
# keeps these as inline strings) instead of auto-coercing to a number.
# ----------------------------------------------------------------------
function Text-Literal($value) {
    return "'" + $value
}

# ========================================================================
# 1. "ODI Batting" sheet (currently sheet 1 of 2): MATCH_CARD_LINK ->
#    MATCH_CODE, full URL -> bare match code text.
# ========================================================================
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2  = "3398"; 3  = "3400"; 4  = "3402"; 5  = "3404"; 6  = "3405";
    7  = "3427"; 8  = "3429"; 9  = "3430"; 10 = "3474"; 11 = "3475";
    12 = "3478"; 13 = "3506"; 14 = "3510"; 15 = "3513"; 16 = "3520";
    17 = "3521"; 18 = "3525"; 19 = "3527"; 20 = "3531"; 21 = "3532";
    22 = "3533"; 23 = "3535"; 24 = "3538"; 25 = "3580"; 26 = "3581";
    27 = "3583"; 28 = "3593"; 29 = "3596"; 30 = "3597"; 31 = "3598";
    32 = "3788"; 33 = "3793"; 34 = "3852"; 35 = "3853"; 36 = "3855";
    37 = "3892"; 38 = "3893"; 39 = "3898"; 40 = "3900"; 41 = "3905";
    42 = "3907"; 43 = "3909"; 44 = "3939"; 45 = "3943"; 46 = "3960";
    47 = "3961"; 48 = "3963"; 49 = "3964"
}

foreach ($r in $battingCodes.Keys) {
    $batting.Cells.Item($r, 4).Value = Text-Literal $battingCodes[$r]
}

# ========================================================================
# 2. "ODI Bowling" sheet (currently sheet 2 of 2): MATCH_CARD_LINK ->
#    MATCH_CODE, full URL -> bare match code text.
# ========================================================================
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = Text-Literal "3852"

# ========================================================================
# 3. New "Player Info" sheet, inserted as the first tab.
# ========================================================================
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$piHeader = $playerInfo.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.HorizontalAlignment = -4108
$piHeader.VerticalAlignment = -4160
$piHeader.Borders.LineStyle = 1

$playerInfo.Range("A2").Value = Text-Literal "3903"
$playerInfo.Range("B2").Value = "Johnson Charles"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

# ========================================================================
# 4. New "ODI Batting Extra" sheet, appended as the last tab.
# ========================================================================
$lastIdx = $wb.Worksheets.Count
$extra = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIdx))
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

$extraHeader = $extra.Range("A1:F1")
$extraHeader.Font.Bold = $true
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160
$extraHeader.Borders.LineStyle = 1

# row data: match code, batting position (number or $null), 4s, 6s,
# % of total runs, man of the match
$extraRows = @(
    @("3597", 2,      "0", "0", $null,    "NO"),
    @("3598", $null,  $null, $null, $null, "NO"),
    @("3788", $null,  $null, $null, $null, "NO"),
    @("3793", 2,      "0", "0", "1.20%",  "NO"),
    @("3852", $null,  $null, $null, $null, "NO"),
    @("3853", 1,      "7", "4", "38.79%", "NO"),
    @("3855", 1,      "1", "0", "1.94%",  "NO"),
    @("3892", 1,      "4", "0", "16.23%", "NO"),
    @("3893", $null,  $null, $null, $null, "NO"),
    @("3898", 1,      "5", "1", "18.05%", "NO"),
    @("3900", $null,  $null, $null, $null, "NO"),
    @("3905", $null,  $null, $null, $null, "NO"),
    @("3907", 1,      "1", "0", "1.40%",  "NO"),
    @("3909", 2,      "4", "1", "21.23%", "NO"),
    @("3939", 1,      "1", "1", "11.43%", "NO"),
    @("3943", 1,      "0", "0", "0.72%",  "NO"),
    @("3960", 1,      "0", "0", "0.88%",  "NO"),
    @("3961", $null,  $null, $null, $null, "NO"),
    @("3963", 1,      "4", "1", "7.90%",  "NO"),
    @("3964", $null,  $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = Text-Literal $row[0]
    if ($null -ne $row[1]) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $extra.Cells.Item($r, 3).Value = Text-Literal $row[2]
    }
    if ($null -ne $row[3]) {
        $extra.Cells.Item($r, 4).Value = Text-Literal $row[3]
    }
    if ($null -ne $row[4]) {
        $extra.Cells.Item($r, 5).Value = Text-Literal $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
